$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H4").Value = $true
$ws.Range("J4").Value = $true
$ws.Range("H5").Value = $true
$ws.Range("J5").Value = $true
$ws.Range("H11").Value = $true
$ws.Range("J11").Value = $true

$ws.Range("G17").Select()
